$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rotation angle column (E) for rows 14-17 from 0 to -90
$ws.Range("E14").Value = -90
$ws.Range("E15").Value = -90
$ws.Range("E16").Value = -90
$ws.Range("E17").Value = -90

# Update the active selection to match the saved view state
$ws.Range("E27").Select()
